$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.50235
$ws.Range("H2").Value = 97.50704999999999
$ws.Range("I2").Value = 0.004318312013857221
$ws.Range("J2").Value = 0.004318312013857221
$ws.Range("M2").Value = 17.08155333333333
$ws.Range("N2").Value = 51.24466
$ws.Range("O2").Value = 0.3501540759902865
$ws.Range("P2").Value = 0.3501540759902865
$ws.Range("Q2").Value = 555.1906249836666
$ws.Range("R2").Value = 4996.715624852999
$ws.Range("S2").Value = 0.001512074553049929
$ws.Range("T2").Value = 0.001512074553049929
$ws.Range("G3").Value = 32.50235
$ws.Range("H3").Value = 97.50704999999999
$ws.Range("I3").Value = 0.004318312013857221
$ws.Range("J3").Value = 0.004318312013857221
$ws.Range("O3").Value = 0.2142771237573249
$ws.Range("P3").Value = 0.2142771237573249
$ws.Range("Q3").Value = 339.74943722155
$ws.Range("R3").Value = 3057.74493499395
$ws.Range("S3").Value = 0.0009253154778160267
$ws.Range("T3").Value = 0.0009253154778160267
$ws.Range("G4").Value = 32.50235
$ws.Range("H4").Value = 97.50704999999999
$ws.Range("I4").Value = 0.004318312013857221
$ws.Range("J4").Value = 0.004318312013857221
$ws.Range("M4").Value = 8.398122666666666
$ws.Range("N4").Value = 25.194368
$ws.Range("O4").Value = 0.1721527793764119
$ws.Range("P4").Value = 0.1721527793764119
$ws.Range("Q4").Value = 272.9587222549333
$ws.Range("R4").Value = 2456.6285002944
$ws.Range("S4").Value = 0.0007434094154000713
$ws.Range("T4").Value = 0.0007434094154000713
$ws.Range("G5").Value = 32.50235
$ws.Range("H5").Value = 97.50704999999999
$ws.Range("I5").Value = 0.004318312013857221
$ws.Range("J5").Value = 0.004318312013857221
$ws.Range("M5").Value = 4.514486333333333
$ws.Range("N5").Value = 13.543459
$ws.Range("O5").Value = 0.09254227409953211
$ws.Range("P5").Value = 0.09254227409953213
$ws.Range("Q5").Value = 146.7314148762167
$ws.Range("R5").Value = 1320.58273388595
$ws.Range("S5").Value = 0.0003996264140336775
$ws.Range("T5").Value = 0.0003996264140336776
$ws.Range("G6").Value = 32.50235
$ws.Range("H6").Value = 97.50704999999999
$ws.Range("I6").Value = 0.004318312013857221
$ws.Range("J6").Value = 0.004318312013857221
$ws.Range("M6").Value = 8.335727666666667
$ws.Range("N6").Value = 25.007183
$ws.Range("O6").Value = 0.1708737467764446
$ws.Range("P6").Value = 0.1708737467764446
$ws.Range("Q6").Value = 270.9307381266834
$ws.Range("R6").Value = 2438.37664314015
$ws.Range("S6").Value = 0.0007378861535575175
$ws.Range("T6").Value = 0.0007378861535575175
$ws.Range("I7").Value = 0.006762540683959845
$ws.Range("J7").Value = 0.006762540683959845
$ws.Range("M7").Value = 17.08155333333333
$ws.Range("N7").Value = 51.24466
$ws.Range("O7").Value = 0.3501540759902865
$ws.Range("P7").Value = 0.3501540759902865
$ws.Range("Q7").Value = 869.4367560188243
$ws.Range("R7").Value = 7824.930804169419
$ws.Range("S7").Value = 0.00236793118453868
$ws.Range("T7").Value = 0.00236793118453868
$ws.Range("I8").Value = 0.006762540683959845
$ws.Range("J8").Value = 0.006762540683959845
$ws.Range("O8").Value = 0.2142771237573249
$ws.Range("P8").Value = 0.2142771237573249
$ws.Range("S8").Value = 0.001449057767050808
$ws.Range("T8").Value = 0.001449057767050808
$ws.Range("I9").Value = 0.006762540683959845
$ws.Range("J9").Value = 0.006762540683959845
$ws.Range("M9").Value = 8.398122666666666
$ws.Range("N9").Value = 25.194368
$ws.Range("O9").Value = 0.1721527793764119
$ws.Range("P9").Value = 0.1721527793764119
$ws.Range("Q9").Value = 427.4574089059128
$ws.Range("R9").Value = 3847.116680153215
$ws.Range("S9").Value = 0.001164190174389749
$ws.Range("T9").Value = 0.001164190174389749
$ws.Range("I10").Value = 0.006762540683959845
$ws.Range("J10").Value = 0.006762540683959845
$ws.Range("M10").Value = 4.514486333333333
$ws.Range("N10").Value = 13.543459
$ws.Range("O10").Value = 0.09254227409953211
$ws.Range("P10").Value = 0.09254227409953213
$ws.Range("Q10").Value = 229.7835727319481
$ws.Range("R10").Value = 2068.052154587533
$ws.Range("S10").Value = 0.0006258208935842494
$ws.Range("T10").Value = 0.0006258208935842495
$ws.Range("I11").Value = 0.006762540683959845
$ws.Range("J11").Value = 0.006762540683959845
$ws.Range("M11").Value = 8.335727666666667
$ws.Range("N11").Value = 25.007183
$ws.Range("O11").Value = 0.1708737467764446
$ws.Range("P11").Value = 0.1708737467764446
$ws.Range("Q11").Value = 424.2815556721246
$ws.Range("R11").Value = 3818.534001049121
$ws.Range("S11").Value = 0.001155540664396359
$ws.Range("T11").Value = 0.001155540664396359
$ws.Range("G12").Value = 3274.382486666667
$ws.Range("H12").Value = 9823.14746
$ws.Range("I12").Value = 0.4350394734576531
$ws.Range("J12").Value = 0.435039473457653
$ws.Range("M12").Value = 17.08155333333333
$ws.Range("N12").Value = 51.24466
$ws.Range("O12").Value = 0.3501540759902865
$ws.Range("P12").Value = 0.3501540759902865
$ws.Range("Q12").Value = 55931.53907972929
$ws.Range("R12").Value = 503383.8517175636
$ws.Range("S12").Value = 0.1523308448478653
$ws.Range("T12").Value = 0.1523308448478652
$ws.Range("G13").Value = 3274.382486666667
$ws.Range("H13").Value = 9823.14746
$ws.Range("I13").Value = 0.4350394734576531
$ws.Range("J13").Value = 0.435039473457653
$ws.Range("O13").Value = 0.2142771237573249
$ws.Range("P13").Value = 0.2142771237573249
$ws.Range("Q13").Value = 34227.3591630482
$ws.Range("R13").Value = 308046.2324674337
$ws.Range("S13").Value = 0.09321900709340698
$ws.Range("T13").Value = 0.09321900709340697
$ws.Range("G14").Value = 3274.382486666667
$ws.Range("H14").Value = 9823.14746
$ws.Range("I14").Value = 0.4350394734576531
$ws.Range("J14").Value = 0.435039473457653
$ws.Range("M14").Value = 8.398122666666666
$ws.Range("N14").Value = 25.194368
$ws.Range("O14").Value = 0.1721527793764119
$ws.Range("P14").Value = 0.1721527793764119
$ws.Range("Q14").Value = 27498.6657806117
$ws.Range("R14").Value = 247487.9920255053
$ws.Range("S14").Value = 0.07489325449418577
$ws.Range("T14").Value = 0.07489325449418575
$ws.Range("G15").Value = 3274.382486666667
$ws.Range("H15").Value = 9823.14746
$ws.Range("I15").Value = 0.4350394734576531
$ws.Range("J15").Value = 0.435039473457653
$ws.Range("M15").Value = 4.514486333333333
$ws.Range("N15").Value = 13.543459
$ws.Range("O15").Value = 0.09254227409953211
$ws.Range("P15").Value = 0.09254227409953213
$ws.Range("Q15").Value = 14782.15498616268
$ws.Range("R15").Value = 133039.3948754641
$ws.Range("S15").Value = 0.04025954219683426
$ws.Range("T15").Value = 0.04025954219683426
$ws.Range("G16").Value = 3274.382486666667
$ws.Range("H16").Value = 9823.14746
$ws.Range("I16").Value = 0.4350394734576531
$ws.Range("J16").Value = 0.435039473457653
$ws.Range("M16").Value = 8.335727666666667
$ws.Range("N16").Value = 25.007183
$ws.Range("O16").Value = 0.1708737467764446
$ws.Range("P16").Value = 0.1708737467764446
$ws.Range("Q16").Value = 27294.36068535614
$ws.Range("R16").Value = 245649.2461682052
$ws.Range("S16").Value = 0.07433682482536082
$ws.Range("T16").Value = 0.0743368248253608
$ws.Range("G17").Value = 7.278837333333333
$ws.Range("H17").Value = 21.836512
$ws.Range("I17").Value = 0.0009670774791190726
$ws.Range("J17").Value = 0.0009670774791190726
$ws.Range("M17").Value = 17.08155333333333
$ws.Range("N17").Value = 51.24466
$ws.Range("O17").Value = 0.3501540759902865
$ws.Range("P17").Value = 0.3501540759902865
$ws.Range("Q17").Value = 124.3338481139911
$ws.Range("R17").Value = 1119.00463302592
$ws.Range("S17").Value = 0.0003386261211119545
$ws.Range("T17").Value = 0.0003386261211119545
$ws.Range("G18").Value = 7.278837333333333
$ws.Range("H18").Value = 21.836512
$ws.Range("I18").Value = 0.0009670774791190726
$ws.Range("J18").Value = 0.0009670774791190726
$ws.Range("O18").Value = 0.2142771237573249
$ws.Range("P18").Value = 0.2142771237573249
$ws.Range("Q18").Value = 76.08621800045866
$ws.Range("R18").Value = 684.7759620041279
$ws.Range("S18").Value = 0.0002072225806761193
$ws.Range("T18").Value = 0.0002072225806761193
$ws.Range("G19").Value = 7.278837333333333
$ws.Range("H19").Value = 21.836512
$ws.Range("I19").Value = 0.0009670774791190726
$ws.Range("J19").Value = 0.0009670774791190726
$ws.Range("M19").Value = 8.398122666666666
$ws.Range("N19").Value = 25.194368
$ws.Range("O19").Value = 0.1721527793764119
$ws.Range("P19").Value = 0.1721527793764119
$ws.Range("Q19").Value = 61.12856879604621
$ws.Range("R19").Value = 550.1571191644159
$ws.Range("S19").Value = 0.0001664850759026823
$ws.Range("T19").Value = 0.0001664850759026823
$ws.Range("G20").Value = 7.278837333333333
$ws.Range("H20").Value = 21.836512
$ws.Range("I20").Value = 0.0009670774791190726
$ws.Range("J20").Value = 0.0009670774791190726
$ws.Range("M20").Value = 4.514486333333333
$ws.Range("N20").Value = 13.543459
$ws.Range("O20").Value = 0.09254227409953211
$ws.Range("P20").Value = 0.09254227409953213
$ws.Range("Q20").Value = 32.86021166388977
$ws.Range("R20").Value = 295.741904975008
$ws.Range("S20").Value = 0.00008949554914812176
$ws.Range("T20").Value = 0.00008949554914812177
$ws.Range("G21").Value = 7.278837333333333
$ws.Range("H21").Value = 21.836512
$ws.Range("I21").Value = 0.0009670774791190726
$ws.Range("J21").Value = 0.0009670774791190726
$ws.Range("M21").Value = 8.335727666666667
$ws.Range("N21").Value = 25.007183
$ws.Range("O21").Value = 0.1708737467764446
$ws.Range("P21").Value = 0.1708737467764446
$ws.Range("Q21").Value = 60.67440574063288
$ws.Range("R21").Value = 546.069651665696
$ws.Range("S21").Value = 0.0001652481522801948
$ws.Range("T21").Value = 0.0001652481522801948
$ws.Range("G22").Value = 4161.570231333333
$ws.Range("H22").Value = 12484.710694
$ws.Range("I22").Value = 0.5529125963654108
$ws.Range("J22").Value = 0.5529125963654108
$ws.Range("M22").Value = 17.08155333333333
$ws.Range("N22").Value = 51.24466
$ws.Range("O22").Value = 0.3501540759902865
$ws.Range("P22").Value = 0.3501540759902865
$ws.Range("Q22").Value = 71086.08385693266
$ws.Range("R22").Value = 639774.754712394
$ws.Range("S22").Value = 0.1936045992837206
$ws.Range("T22").Value = 0.1936045992837206
$ws.Range("G23").Value = 4161.570231333333
$ws.Range("H23").Value = 12484.710694
$ws.Range("I23").Value = 0.5529125963654108
$ws.Range("J23").Value = 0.5529125963654108
$ws.Range("O23").Value = 0.2142771237573249
$ws.Range("P23").Value = 0.2142771237573249
$ws.Range("Q23").Value = 43501.19742275422
$ws.Range("R23").Value = 391510.776804788
$ws.Range("S23").Value = 0.118476520838375
$ws.Range("T23").Value = 0.118476520838375
$ws.Range("G24").Value = 4161.570231333333
$ws.Range("H24").Value = 12484.710694
$ws.Range("I24").Value = 0.5529125963654108
$ws.Range("J24").Value = 0.5529125963654108
$ws.Range("M24").Value = 8.398122666666666
$ws.Range("N24").Value = 25.194368
$ws.Range("O24").Value = 0.1721527793764119
$ws.Range("P24").Value = 0.1721527793764119
$ws.Range("Q24").Value = 34949.3772886857
$ws.Range("R24").Value = 314544.3955981713
$ws.Range("S24").Value = 0.09518544021653366
$ws.Range("T24").Value = 0.09518544021653366
$ws.Range("G25").Value = 4161.570231333333
$ws.Range("H25").Value = 12484.710694
$ws.Range("I25").Value = 0.5529125963654108
$ws.Range("J25").Value = 0.5529125963654108
$ws.Range("M25").Value = 4.514486333333333
$ws.Range("N25").Value = 13.543459
$ws.Range("O25").Value = 0.09254227409953211
$ws.Range("P25").Value = 0.09254227409953213
$ws.Range("Q25").Value = 18787.35193456117
$ws.Range("R25").Value = 169086.1674110506
$ws.Range("S25").Value = 0.05116778904593181
$ws.Range("T25").Value = 0.05116778904593182
$ws.Range("G26").Value = 4161.570231333333
$ws.Range("H26").Value = 12484.710694
$ws.Range("I26").Value = 0.5529125963654108
$ws.Range("J26").Value = 0.5529125963654108
$ws.Range("M26").Value = 8.335727666666667
$ws.Range("N26").Value = 25.007183
$ws.Range("O26").Value = 0.1708737467764446
$ws.Range("P26").Value = 0.1708737467764446
$ws.Range("Q26").Value = 34689.71611410167
$ws.Range("R26").Value = 312207.445026915
$ws.Range("S26").Value = 0.09447824698084974
$ws.Range("T26").Value = 0.09447824698084974
